$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: cells changing from TEXT style to NUMERIC style ---
# (copy style from untouched numeric donors G22 [style 15, #,##0] / H22 [style 16, pct fmt], then set value)
$ws.Range("G22").Copy($ws.Range("C14"))
$ws.Range("G22").Copy($ws.Range("F14"))
$ws.Range("G22").Copy($ws.Range("D26"))
$ws.Range("H22").Copy($ws.Range("E26"))
$ws.Range("G22").Copy($ws.Range("C28"))
$ws.Range("G22").Copy($ws.Range("F28"))
$ws.Range("G22").Copy($ws.Range("C29"))
$ws.Range("G22").Copy($ws.Range("F29"))

# --- Step 2: cells changing from NUMERIC style to TEXT style ---
# (copy style+value from untouched text donors C22 ["0", style 14] / E22 ["***.*", style 14])
$ws.Range("C22").Copy($ws.Range("D27"))
$ws.Range("E22").Copy($ws.Range("E27"))
$ws.Range("C22").Copy($ws.Range("F27"))
$ws.Range("C22").Copy($ws.Range("D28"))
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("C22").Copy($ws.Range("D29"))
$ws.Range("E22").Copy($ws.Range("E29"))
$ws.Range("C22").Copy($ws.Range("D30"))
$ws.Range("E22").Copy($ws.Range("E30"))

# --- Step 3: set numeric values (covers cells touched in Step 1 plus all other numeric-only updates) ---
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = -20
$ws.Range("L14").Value = -20
$ws.Range("M14").Value = -55.555555555555
$ws.Range("N14").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 63.636363636363
$ws.Range("L15").Value = 5.882352941176
$ws.Range("M15").Value = 38.461538461538
$ws.Range("N15").Value = -25
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = -6.896551724137
$ws.Range("I16").Value = 183
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = 11.585365853658
$ws.Range("L16").Value = 33.576642335766
$ws.Range("M16").Value = -3.684210526315
$ws.Range("N16").Value = -59.602649006622
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 7.142857142857
$ws.Range("I17").Value = 265
$ws.Range("J17").Value = 219
$ws.Range("K17").Value = 21.004566210045
$ws.Range("L17").Value = 43.243243243243
$ws.Range("M17").Value = 54.970760233918
$ws.Range("N17").Value = 28.640776699029
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 157
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 72.527472527472
$ws.Range("L18").Value = 50.961538461538
$ws.Range("M18").Value = -34.309623430962
$ws.Range("N18").Value = -83.780991735537
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 81.818181818181
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 35.416666666666
$ws.Range("I19").Value = 421
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 5.25
$ws.Range("L19").Value = 57.089552238806
$ws.Range("M19").Value = 49.290780141844
$ws.Range("N19").Value = 9.067357512953
$ws.Range("C20").Value = 9
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 45
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 87.5
$ws.Range("I20").Value = 350
$ws.Range("J20").Value = 231
$ws.Range("K20").Value = 51.515151515151
$ws.Range("L20").Value = 127.272727272727
$ws.Range("M20").Value = 128.758169934641
$ws.Range("N20").Value = -71.498371335504
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = 35.294117647058
$ws.Range("F21").Value = 185
$ws.Range("H21").Value = 31.205673758865
$ws.Range("I21").Value = 1398
$ws.Range("J21").Value = 1121
$ws.Range("K21").Value = 24.710080285459
$ws.Range("L21").Value = 60.689655172413
$ws.Range("M21").Value = 32.261116367076
$ws.Range("N21").Value = -57.286892758936
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = 100
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 67
$ws.Range("K23").Value = 32.835820895522
$ws.Range("L23").Value = 53.448275862069
$ws.Range("M23").Value = 97.777777777777
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 12
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 3.389830508474
$ws.Range("I24").Value = 1047
$ws.Range("J24").Value = 872
$ws.Range("K24").Value = 20.068807339449
$ws.Range("L24").Value = 59.118541033434
$ws.Range("M24").Value = 72.487644151565
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -61.538461538461
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -35
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 367
$ws.Range("K25").Value = -3.269754768392
$ws.Range("L25").Value = 48.535564853556
$ws.Range("M25").Value = -12.990196078431
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -3.703703703703
$ws.Range("L26").Value = 13.043478260869
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = 104.347826086957
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = -52.173913043478
$ws.Range("L28").Value = 10
$ws.Range("M28").Value = -45
$ws.Range("N28").Value = -56
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 10
$ws.Range("K29").Value = -44.444444444444
$ws.Range("L29").Value = 11.111111111111
$ws.Range("M29").Value = -37.5
$ws.Range("N29").Value = -52.380952380952

# --- Step 4: shared text updates (title volume number + report week range) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"
